$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.368.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.478.42"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.917.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.291.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.23"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.482.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "322.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.409"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0754"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.22"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.60%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.800"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "277.85"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0910"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.743.16"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -1.50%  "
